# Scheduled runner update: refresh cached market-board price/profit figures
# (currentAveragePrice*, LevePrice*, LeveProfit*) across the ALC/ARM/BSM/CRP/
# GSM/LTW/WVR leve-profit sheets with the latest pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 7148.9
$ws.Range("I43").Value = 1996.3334
$ws.Range("K43").Value = 1996.3334
$ws.Range("M43").Value = -1927.3334

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 1111381.9
$ws.Range("J96").Value = 278.33334
$ws.Range("L96").Value = 835.0000200000001
$ws.Range("N96").Value = -3581.00002

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 3890.3
$ws.Range("I100").Value = 3655.889
$ws.Range("K100").Value = 3655.889
$ws.Range("M100").Value = -3114.889

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 12957.5
$ws.Range("I137").Value = 22449
$ws.Range("K137").Value = 67347
$ws.Range("M137").Value = -64797

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3454.6206
$ws.Range("J138").Value = 3663.2666
$ws.Range("L138").Value = 10989.7998
$ws.Range("N138").Value = -21269.7998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 658.6923
$ws.Range("I97").Value = 151.44444
$ws.Range("J97").Value = 1800
$ws.Range("K97").Value = 151.44444
$ws.Range("L97").Value = 1800
$ws.Range("M97").Value = 344.55556
$ws.Range("N97").Value = -2792

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H119").Value = 105249
$ws.Range("J119").Value = 105249
$ws.Range("L119").Value = 105249
$ws.Range("N119").Value = -114925

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 6538906.5
$ws.Range("I122").Value = 7939498.5
$ws.Range("J122").Value = 2809
$ws.Range("K122").Value = 23818495.5
$ws.Range("L122").Value = 8427
$ws.Range("M122").Value = -23816045.5
$ws.Range("N122").Value = -13327

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H124").Value = 11750
$ws.Range("J124").Value = 11750
$ws.Range("L124").Value = 11750
$ws.Range("N124").Value = -21570

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H125").Value = 62857.5
$ws.Range("J125").Value = 62857.5
$ws.Range("L125").Value = 62857.5
$ws.Range("N125").Value = -72697.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2905.2856
$ws.Range("I132").Value = 2937.7778
$ws.Range("J132").Value = 2880.9167
$ws.Range("K132").Value = 8813.3334
$ws.Range("L132").Value = 8642.750100000001
$ws.Range("M132").Value = -6283.3334
$ws.Range("N132").Value = -13702.7501

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H23").Value = 3733.3333
$ws.Range("I23").Value = 3100
$ws.Range("J23").Value = 5000
$ws.Range("K23").Value = 3100
$ws.Range("L23").Value = 5000
$ws.Range("M23").Value = -2817
$ws.Range("N23").Value = -5566

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 4482.5
$ws.Range("I107").Value = 1062.3182
$ws.Range("K107").Value = 1062.3182
$ws.Range("M107").Value = 857.6818000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3930.1428
$ws.Range("J134").Value = 3133
$ws.Range("L134").Value = 9399
$ws.Range("N134").Value = -14469

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1084.1177
$ws.Range("J16").Value = 1383.8
$ws.Range("L16").Value = 1383.8
$ws.Range("N16").Value = -1957.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("M27").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3827.2778
$ws.Range("I31").Value = 1762.7693
$ws.Range("J31").Value = 9195
$ws.Range("K31").Value = 1762.7693
$ws.Range("L31").Value = 9195
$ws.Range("M31").Value = -1467.7693
$ws.Range("N31").Value = -9785

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3827.2778
$ws.Range("I34").Value = 1762.7693
$ws.Range("J34").Value = 9195
$ws.Range("K34").Value = 1762.7693
$ws.Range("L34").Value = 9195
$ws.Range("M34").Value = -1560.7693
$ws.Range("N34").Value = -9599

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 2583.6
$ws.Range("I105").Value = 1844.5
$ws.Range("J105").Value = 3428.2856
$ws.Range("K105").Value = 1844.5
$ws.Range("L105").Value = 3428.2856
$ws.Range("M105").Value = -97.5
$ws.Range("N105").Value = -6922.2856

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 1084.1177
$ws.Range("J113").Value = 1383.8
$ws.Range("L113").Value = 1383.8
$ws.Range("N113").Value = -5723.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1759.25
$ws.Range("I132").Value = 1832.8182
$ws.Range("J132").Value = 950
$ws.Range("K132").Value = 5498.4546
$ws.Range("L132").Value = 2850
$ws.Range("M132").Value = -2968.4546
$ws.Range("N132").Value = -7910

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8780718
$ws.Range("I70").Value = 55562840
$ws.Range("J70").Value = 9069.9375
$ws.Range("K70").Value = 55562840
$ws.Range("L70").Value = 9069.9375
$ws.Range("M70").Value = -55562570
$ws.Range("N70").Value = -9609.9375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 8780718
$ws.Range("I73").Value = 55562840
$ws.Range("J73").Value = 9069.9375
$ws.Range("K73").Value = 55562840
$ws.Range("L73").Value = 9069.9375
$ws.Range("M73").Value = -55561904
$ws.Range("N73").Value = -10941.9375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 4488.6895
$ws.Range("I97").Value = 1202.6957
$ws.Range("K97").Value = 1202.6957
$ws.Range("M97").Value = -706.6957

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 289.27274
$ws.Range("I107").Value = 213.28572
$ws.Range("J107").Value = 422.25
$ws.Range("K107").Value = 213.28572
$ws.Range("L107").Value = 422.25
$ws.Range("M107").Value = 1706.71428
$ws.Range("N107").Value = -4262.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2614.5173
$ws.Range("I132").Value = 2266.1738
$ws.Range("K132").Value = 6798.5214
$ws.Range("M132").Value = -4268.5214

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3724.1875
$ws.Range("I22").Value = 3113.3572
$ws.Range("K22").Value = 3113.3572
$ws.Range("M22").Value = -2818.3572

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 3724.1875
$ws.Range("I27").Value = 3113.3572
$ws.Range("K27").Value = 3113.3572
$ws.Range("M27").Value = -3006.3572

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 427972.62
$ws.Range("I82").Value = 756751.5600000001
$ws.Range("K82").Value = 756751.5600000001
$ws.Range("M82").Value = -756390.5600000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 427972.62
$ws.Range("I85").Value = 756751.5600000001
$ws.Range("K85").Value = 756751.5600000001
$ws.Range("M85").Value = -755503.5600000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5899.2666
$ws.Range("I132").Value = 6063
$ws.Range("J132").Value = 5449
$ws.Range("K132").Value = 18189
$ws.Range("L132").Value = 16347
$ws.Range("M132").Value = -15659
$ws.Range("N132").Value = -21407

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 8077.222
$ws.Range("I96").Value = 8077.222
$ws.Range("K96").Value = 8077.222
$ws.Range("M96").Value = -6704.222

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1067.375
$ws.Range("I107").Value = 1236.9231
$ws.Range("J107").Value = 332.66666
$ws.Range("K107").Value = 3710.7693
$ws.Range("L107").Value = 997.9999799999999
$ws.Range("M107").Value = -1790.7693
$ws.Range("N107").Value = -4837.99998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2946.65
$ws.Range("I126").Value = 2852.5
$ws.Range("J126").Value = 3166.3333
$ws.Range("K126").Value = 8557.5
$ws.Range("L126").Value = 9498.999899999999
$ws.Range("M126").Value = -6087.5
$ws.Range("N126").Value = -14438.9999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2449.7727
$ws.Range("I136").Value = 2219.75
$ws.Range("K136").Value = 6659.25
$ws.Range("M136").Value = -4109.25
